# Add the "weather" localization rows (Sunny / Partly Sunny / Mostly Cloudy /
# Light Rain / Rain) to the language sheet, right after the existing
# day_sunday / Sunday row (row 44) -> rows 45-49, columns A (key) / B (value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("weather_sunny", "Sunny"),
    @("weather_partly_sunny", "Partly Sunny"),
    @("weather_mostly_cloudy", "Mostly Cloudy"),
    @("weather_light_rain", "Light Rain"),
    @("weather_rain", "Rain")
)

$row = 45
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]

    $cellB = $ws.Cells.Item($row, 2)
    $cellB.Value = $pair[1]
    $cellB.WrapText = $true

    $row = $row + 1
}

# Mirror the author's final selection / scroll position: active cell A49,
# scrolled so row 7 is at the top of the viewport.
$ws.Range("A49").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1

$wb.Save()
